$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 1.0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text: fix "usuario" -> "usuário" and add trailing period.
# This precondition text is repeated for each of the four test cases (TC1-TC4).
$newPrecondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value = $newPrecondition
$ws.Range("B16").Value = $newPrecondition
$ws.Range("B24").Value = $newPrecondition
$ws.Range("B32").Value = $newPrecondition

# TC1 step 2 expected result: add trailing period
$ws.Range("D11").Value = "SYSTEM Apresenta a tela de Registrar Liquidações."

# TC2 step 2 expected result: remove duplicate "o nome"
$ws.Range("D19").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# TC3 step 2 expected result: add trailing period
$ws.Range("D27").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."

# TC4 step 2 expected result: accent fixes (numero -> número, diaria -> diária, diarias -> diárias)
$ws.Range("D35").Value = "SYSTEM Exibe a lista de diárias (solicitações) aptas para pagamento ordenado pelo número da diária em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de liquidação (após registrar o empenho)."
